# Auto update Excel log
# Appends newly-logged sensor readings to the PIR, Humidity, and Temperature sheets.

function Set-TextValue($ws, $r, $c, $val) {
    # Force number format to Text first so Excel does not auto-convert
    # date-looking / percentage-looking strings into numeric values.
    $ws.Cells.Item($r, $c).NumberFormat = "@"
    $ws.Cells.Item($r, $c).Value = $val
}

function Add-SensorRow($ws, $r, $date, $timestamp, $hour, $location, $value, $status) {
    Set-TextValue $ws $r 1 $date
    Set-TextValue $ws $r 2 $timestamp
    Set-TextValue $ws $r 3 $hour
    Set-TextValue $ws $r 4 $location
    Set-TextValue $ws $r 5 $value
    Set-TextValue $ws $r 6 $status
}

$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 21-33 (No Motion / Inactive) ---
$wsPIR = $wb.Worksheets.Item("PIR")
Add-SensorRow $wsPIR 21 "2026-02-01" "18:02:15" "18:00" "Bathroom" "No Motion" "Inactive"
Add-SensorRow $wsPIR 22 "2026-02-01" "18:02:16" "18:00" "Bathroom" "No Motion" "Inactive"
Add-SensorRow $wsPIR 23 "2026-02-01" "18:02:21" "18:00" "Bathroom" "No Motion" "Inactive"
Add-SensorRow $wsPIR 24 "2026-02-01" "18:02:26" "18:00" "Bathroom" "No Motion" "Inactive"
Add-SensorRow $wsPIR 25 "2026-02-01" "18:02:31" "18:00" "Bathroom" "No Motion" "Inactive"
Add-SensorRow $wsPIR 26 "2026-02-01" "18:02:36" "18:00" "Bathroom" "No Motion" "Inactive"
Add-SensorRow $wsPIR 27 "2026-02-01" "18:02:41" "18:00" "Bathroom" "No Motion" "Inactive"
Add-SensorRow $wsPIR 28 "2026-02-01" "18:02:46" "18:00" "Bathroom" "No Motion" "Inactive"
Add-SensorRow $wsPIR 29 "2026-02-01" "18:02:51" "18:00" "Bathroom" "No Motion" "Inactive"
Add-SensorRow $wsPIR 30 "2026-02-01" "18:02:56" "18:00" "Bathroom" "No Motion" "Inactive"
Add-SensorRow $wsPIR 31 "2026-02-01" "18:03:01" "18:00" "Bathroom" "No Motion" "Inactive"
Add-SensorRow $wsPIR 32 "2026-02-01" "18:03:06" "18:00" "Bathroom" "No Motion" "Inactive"
Add-SensorRow $wsPIR 33 "2026-02-01" "18:03:11" "18:00" "Bathroom" "No Motion" "Inactive"

# --- Humidity sheet: append rows 19-31 (Active) ---
$wsHumidity = $wb.Worksheets.Item("Humidity")
Add-SensorRow $wsHumidity 19 "2026-02-01" "18:02:16" "18:00" "Bathroom" "80.7%" "Active"
Add-SensorRow $wsHumidity 20 "2026-02-01" "18:02:16" "18:00" "Bathroom" "81.6%" "Active"
Add-SensorRow $wsHumidity 21 "2026-02-01" "18:02:21" "18:00" "Bathroom" "80.6%" "Active"
Add-SensorRow $wsHumidity 22 "2026-02-01" "18:02:26" "18:00" "Bathroom" "81.5%" "Active"
Add-SensorRow $wsHumidity 23 "2026-02-01" "18:02:31" "18:00" "Bathroom" "80.5%" "Active"
Add-SensorRow $wsHumidity 24 "2026-02-01" "18:02:36" "18:00" "Bathroom" "81.4%" "Active"
Add-SensorRow $wsHumidity 25 "2026-02-01" "18:02:41" "18:00" "Bathroom" "80.4%" "Active"
Add-SensorRow $wsHumidity 26 "2026-02-01" "18:02:46" "18:00" "Bathroom" "80.0%" "Active"
Add-SensorRow $wsHumidity 27 "2026-02-01" "18:02:51" "18:00" "Bathroom" "80.3%" "Active"
Add-SensorRow $wsHumidity 28 "2026-02-01" "18:02:57" "18:00" "Bathroom" "81.2%" "Active"
Add-SensorRow $wsHumidity 29 "2026-02-01" "18:03:01" "18:00" "Bathroom" "80.2%" "Active"
Add-SensorRow $wsHumidity 30 "2026-02-01" "18:03:06" "18:00" "Bathroom" "81.1%" "Active"
Add-SensorRow $wsHumidity 31 "2026-02-01" "18:03:11" "18:00" "Bathroom" "80.1%" "Active"

# --- Temperature sheet: append rows 19-31 (Active) ---
$wsTemperature = $wb.Worksheets.Item("Temperature")
Add-SensorRow $wsTemperature 19 "2026-02-01" "18:02:16" "18:00" "Bathroom" "28.7C" "Active"
Add-SensorRow $wsTemperature 20 "2026-02-01" "18:02:17" "18:00" "Bathroom" "28.8C" "Active"
Add-SensorRow $wsTemperature 21 "2026-02-01" "18:02:21" "18:00" "Bathroom" "28.8C" "Active"
Add-SensorRow $wsTemperature 22 "2026-02-01" "18:02:26" "18:00" "Bathroom" "28.8C" "Active"
Add-SensorRow $wsTemperature 23 "2026-02-01" "18:02:31" "18:00" "Bathroom" "28.8C" "Active"
Add-SensorRow $wsTemperature 24 "2026-02-01" "18:02:36" "18:00" "Bathroom" "28.8C" "Active"
Add-SensorRow $wsTemperature 25 "2026-02-01" "18:02:41" "18:00" "Bathroom" "28.8C" "Active"
Add-SensorRow $wsTemperature 26 "2026-02-01" "18:02:46" "18:00" "Bathroom" "28.8C" "Active"
Add-SensorRow $wsTemperature 27 "2026-02-01" "18:02:51" "18:00" "Bathroom" "28.8C" "Active"
Add-SensorRow $wsTemperature 28 "2026-02-01" "18:02:57" "18:00" "Bathroom" "28.8C" "Active"
Add-SensorRow $wsTemperature 29 "2026-02-01" "18:03:01" "18:00" "Bathroom" "28.8C" "Active"
Add-SensorRow $wsTemperature 30 "2026-02-01" "18:03:06" "18:00" "Bathroom" "28.8C" "Active"
Add-SensorRow $wsTemperature 31 "2026-02-01" "18:03:12" "18:00" "Bathroom" "28.8C" "Active"

